$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.061.01"
$ws.Range("D3").Value = "1.668.10"
$ws.Range("E3").Value = "  -1.57%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.84"
$ws.Range("E5").Value = "  -1.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5113"
$ws.Range("E6").Value = "  +0.48%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2656"
$ws.Range("E8").Value = "  +0.43%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06407"
$ws.Range("E9").Value = "  +1.90%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.88"
$ws.Range("E10").Value = "  -0.99%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07429"
$ws.Range("E11").Value = "  +0.73%  "
$ws.Range("D12").Value = "1.683.78"
$ws.Range("E12").Value = "  -0.73%  "
$ws.Range("E13").Value = "  -0.29%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5841"
$ws.Range("E14").Value = "  +0.24%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000008543"
$ws.Range("E15").Value = "  +1.87%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.32"
$ws.Range("E16").Value = "  -1.73%  "
$ws.Range("D17").Value = "26.079.51"
$ws.Range("E17").Value = "  -2.03%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.945"
$ws.Range("E18").Value = "  -1.28%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.005"
$ws.Range("E19").Value = "  -0.07%  "
$ws.Range("E20").Value = "  -2.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "190.35"
$ws.Range("E21").Value = "  +2.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.224"
$ws.Range("E22").Value = "  -0.58%  "
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "145.13"
$ws.Range("E24").Value = "  +0.46%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "7.611"
$ws.Range("E25").Value = "  +1.27%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1201"
$ws.Range("E26").Value = "  +3.87%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.60"
$ws.Range("E27").Value = "  -0.41%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06492"
$ws.Range("E28").Value = "  +14.82%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.317"
$ws.Range("E29").Value = "  -1.88%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.316"
$ws.Range("E30").Value = "  -1.53%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.535"
$ws.Range("E31").Value = "  +0.75%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.516"
$ws.Range("E32").Value = "  +1.01%  "
$ws.Range("E33").Value = "  +0.52%  "
$ws.Range("E34").Value = "  -0.19%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6094"
$ws.Range("E35").Value = "  +1.00%  "
$ws.Range("E36").Value = "  +0.18%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.716"
$ws.Range("E37").Value = "  +1.35%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.233"
$ws.Range("E38").Value = "  +6.58%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01601"
$ws.Range("E39").Value = "  -0.70%  "
$ws.Range("D40").Value = "1.083.73"
$ws.Range("E40").Value = "  -0.99%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8661"
$ws.Range("E41").Value = "  +0.90%  "
$ws.Range("E42").Value = "  +0.71%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.63"
$ws.Range("E43").Value = "  +1.05%  "
$ws.Range("D44").Value = "1.816.79"
$ws.Range("E44").Value = "  -1.94%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000113"
$ws.Range("E45").Value = "  +3.63%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.29"
$ws.Range("E46").Value = "  -0.71%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.007"
$ws.Range("E47").Value = "  +0.29%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.098"
$ws.Range("E48").Value = "  -1.22%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05239"
$ws.Range("E49").Value = "  -0.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4288"
$ws.Range("E50").Value = "  -0.88%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.019"
$ws.Range("E51").Value = "  +4.14%  "
